$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Copy the date-formatted style from the last existing date row (G17) so
# the new date cells use the same cell style (numFmtId 14) as the rest
# of the column, rather than Excel's auto date-format heuristics.
$ws.Range("G17").Copy() | Out-Null

# New rows of timesheet data
# Row 18
$ws.Range("G18").Value = 43662
$ws.Range("G18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H18").Value = 0.75
$ws.Range("I18").Value = "Stunden"
$ws.Range("J18").Value = "Dokumentation, Statusbericht"

# Row 19
$ws.Range("G19").Value = 43664
$ws.Range("G19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = "Stunden"
$ws.Range("J19").Value = "Recherche"
$ws.Range("K19").Value = "Klassen und Bibilotheken in IntelliJ IDEA"

# Row 20
$ws.Range("G20").Value = 43666
$ws.Range("G20").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H20").Value = 0.75
$ws.Range("I20").Value = "Stunden"
$ws.Range("J20").Value = "Apache Ant"

# Row 21
$ws.Range("G21").Value = 43667
$ws.Range("G21").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = "Stunden"
$ws.Range("J21").Value = "Programmieren"
$ws.Range("K21").Value = "Erstellen neuer Klassen, Package testui zum Testen"

$excel.CutCopyMode = 0

# Recalculate formulas (B6 = SUM(H:H), B7 depends on B6 and TODAY())
$excel.Calculate()

# Update selection to match the diff
$ws.Range("G22").Select()
